$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the new refGeo function data
$ws.Range("A2").Value = 15
$ws.Range("C2").Value = 7375
$ws.Range("D2").Value = 70740
$ws.Range("E2").Value = 0.7417
$ws.Range("F2").Value = "F02"

# G2 and H2 hold numeric-looking codes that must remain stored as text,
# so format them as text before assigning, then restore the default style
# (the underlying value keeps its text/shared-string type either way).
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "17"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "72752"
$ws.Range("H2").Style = "Normal"
